# Applies updated input data to the InstrumentAttribute sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("InstrumentAttribute")

# Update Notional/Quantity style amounts for rows 23-26 (columns H and J)
$ws.Range("H23").Value = -1000
$ws.Range("J23").Value = -1000

$ws.Range("H24").Value = 500
$ws.Range("J24").Value = 6000

$ws.Range("H25").Value = 1500
$ws.Range("J25").Value = 1500

$ws.Range("H26").Value = 250
$ws.Range("J26").Value = 6000

# Reflect the active selection left in the sheet as J23
$ws.Activate()
$ws.Range("J23").Select()
